$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue "D2" "44.151.82"
Set-TextValue "E2" "  +0.51%  "
Set-TextValue "D3" "2.246.15"
Set-TextValue "E3" "  +0.56%  "
Set-TextValue "E4" "  +0.00%  "
Set-TextValue "D5" "306.47"
Set-TextValue "E5" "  -1.87%  "
Set-TextValue "D6" "95.87"
Set-TextValue "E6" "  -2.76%  "
Set-TextValue "E7" "  +1.09%  "
Set-TextValue "E9" "  -1.22%  "
Set-TextValue "D10" "34.96"
Set-TextValue "E10" "  -3.22%  "
Set-TextValue "D11" "0.0818"
Set-TextValue "E11" "  -0.14%  "
Set-TextValue "E12" "  -1.21%  "
Set-TextValue "D13" "0.105"
Set-TextValue "E13" "  +0.30%  "
Set-TextValue "D14" "2.361.44"
Set-TextValue "E14" "  +4.38%  "
Set-TextValue "D15" "2.587.94"
Set-TextValue "E15" "  +0.50%  "
Set-TextValue "E16" "  -0.23%  "
Set-TextValue "D17" "13.64"
Set-TextValue "E17" "  -3.17%  "
Set-TextValue "D18" "44.059.30"
Set-TextValue "E18" "  +0.59%  "
Set-TextValue "D19" "0.0₃0976"
Set-TextValue "E19" "  +1.70%  "
Set-TextValue "D20" "12.20"
Set-TextValue "E20" "  -5.25%  "
Set-TextValue "E21" "  +1.07%  "
Set-TextValue "E22" "  +0.85%  "
Set-TextValue "D23" "236.88"
Set-TextValue "E23" "  +1.74%  "
Set-TextValue "E24" "  -1.12%  "
Set-TextValue "E25" "  -1.11%  "
Set-TextValue "E26" "  +0.00%  "
Set-TextValue "D27" "9.96"
Set-TextValue "E27" "  -1.56%  "
Set-TextValue "D29" "37.62"
Set-TextValue "E29" "  +2.64%  "
Set-TextValue "D30" "6.01"
Set-TextValue "E30" "  +1.21%  "
Set-TextValue "D31" "20.10"
Set-TextValue "E31" "  +1.08%  "
Set-TextValue "D32" "152.32"
Set-TextValue "E32" "  -3.46%  "
Set-TextValue "E33" "  -2.65%  "
Set-TextValue "D34" "3.32"
Set-TextValue "E34" "  +4.05%  "
Set-TextValue "E35" "  -3.02%  "
Set-TextValue "B36" "Kaspa"
Set-TextValue "C36" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D36" "0.110"
Set-TextValue "E36" "  +0.59%  "
Set-TextValue "B37" "Stellar"
Set-TextValue "C37" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D37" "0.120"
Set-TextValue "E37" "  +2.51%  "
Set-TextValue "E38" "  -6.37%  "
Set-TextValue "B39" "RenderToken"
Set-TextValue "C39" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D39" "3.88"
Set-TextValue "E39" "  -3.80%  "
Set-TextValue "B40" "Celestia"
Set-TextValue "C40" "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue "D40" "14.63"
Set-TextValue "E40" "  -5.97%  "
Set-TextValue "D41" "3.39"
Set-TextValue "E41" "  -5.53%  "
Set-TextValue "E42" "  -2.56%  "
Set-TextValue "E43" "  +0.18%  "
Set-TextValue "D44" "1.738.60"
Set-TextValue "E44" "  +1.63%  "
Set-TextValue "D45" "83.55"
Set-TextValue "E45" "  +4.35%  "
Set-TextValue "E46" "  -1.25%  "
Set-TextValue "D47" "100.43"
Set-TextValue "D48" "4.90"
Set-TextValue "E48" "  -3.48%  "
Set-TextValue "E49" "  +1.84%  "
Set-TextValue "D50" "54.80"
Set-TextValue "E50" "  -2.48%  "
Set-TextValue "D51" "68.33"
Set-TextValue "E51" "  -6.32%  "
